# Use DefaultTenant as default cloud tenant name
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Row 3 holds the CloudTenantName setting: B3 = Value, C3 = Explanation
$ws.Range("B3").Value = "DefaultTenant"
$ws.Range("C3").Value = "Name of the tenant to be used in case of Automation Cloud Orchestrator instances. `nSample value: DefaultTenant."
